# fix(publipostage): Try to solve Excel emoji problem
#
# Replace the 4 "status" emoji markers used in column A with
# non-pictographic equivalents that survive mail-merge / publipostage
# rendering better:
#   📘 (blue book)   -> ⚠️
#   📕 (red book)    -> -3
#   📗 (green book)  -> ✅
#   📙 (orange book) -> +3
#
# "-3" and "+3" look like numbers to Excel's automatic type inference, so
# assigning them naively via .Value2 turns the cell into a numeric cell
# (losing the leading sign / text formatting). We force them to stay text
# by entering them with a leading quote (like typing '-3 into Excel), then
# restore the cell's original style so no extra formatting sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📗" = "✅"
    "📙" = "+3"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($null -ne $old -and $map.ContainsKey($old)) {
        $new = $map[$old]
        $origStyle = $cell.Style
        if ($new -eq "-3" -or $new -eq "+3") {
            # Enter with a leading apostrophe so Excel keeps it as text
            # instead of converting it to a signed number.
            $cell.Formula = "'" + $new
        } else {
            $cell.Value2 = $new
        }
        $cell.Style = $origStyle
    }
}
